# Deviceinventories.xlsx migration edit
# Commit message: "added swagger and migration"
#
# The underlying change: the "fontId" column (column P, the sheet's header
# in P1 held the shared string "fontId") was removed from the worksheet
# entirely. Deleting that column shifts every column at/after Q left by
# one (Q->P ... AQ->AP) and the SQL-building helper formula (previously
# in column S, referencing the now-removed P column via "&P{row}&" and
# including ", fontId" in the generated INSERT statement) moves to column
# R and no longer mentions fontId at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Delete the entire "fontId" column (column P / 16).
#    This shifts every column after it one to the left, so the helper
#    formula column (originally S) becomes R, etc. Any formula that
#    referenced column P (the "&P2&" piece) will temporarily evaluate to
#    #REF! until we rewrite those formulas below.
$ws.Columns.Item(16).Delete()

# 2) Rewrite the SQL-insert helper formulas (now in column R, rows 2-45)
#    so they no longer reference the deleted fontId column or mention it
#    in the generated INSERT statement text.
for ($r = 2; $r -le 45; $r++) {
    $nextRow = $r + 1
    if ($r -eq 45) { $nextRow = 45 }

    $formula = '="INSERT INTO bubbldb.deviceinventories (id, name, deviceTypeId, productId, shortDescription, deviceDescription, materialTypeId, patternId, colorId, productDetails, price, discountPercentage, availability, createdAt, updatedAt) VALUES (''"&A' + $r + '&"'',''"&B' + $r + '&"'',''"&C' + $r + '&"'',''"&D' + $r + '&"'',''"&E' + $r + '&"'',''"&F' + $r + '&"'',''"&G' + $r + '&"'',''"&H' + $r + '&"'',''"&I' + $r + '&"'',''"&J' + $r + '&"'',''"&K' + $r + '&"'',''"&L' + $r + '&"'',''"&M' + $r + '&"'',''"&TEXT(O' + $nextRow + ', "yyyy-mm-dd hh:mm:ss")&"'',''"&TEXT(O' + $nextRow + ', "yyyy-mm-dd hh:mm:ss")&"'');"'

    $ws.Range("R$r").Formula = $formula
}

$excel.Calculate()

# 3) Reflect the viewport/selection state recorded in the saved file
#    (cursor had scrolled down and over before the author saved).
$ws.Range("A32").Select()
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("S51").Select()

$wb.Save()
